$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DecisionTable")

# Row 18: C18 changes from "ACTION" to "CONDITION" (value 7); new D18/E18 = "ACTION" (value 5)
$ws.Range("C18").Value = "CONDITION"
$ws.Range("D18").Value = "ACTION"
$ws.Range("E18").Value = "ACTION"

# Row 19: C19 changes from "Test" to "Code changed 10010018"; D19 = "Test"; E19 = "Code changed 10010018"
$ws.Range("C19").Value = "Code changed 10010018"
$ws.Range("D19").Value = "Test"
$ws.Range("E19").Value = "Code changed 10010018"

# Row 20: C20 changes from "gfndnvbx" to "Code changed 10010018"; D20 = "gfndnvbx"; E20 = "Code changed 10010018"
$ws.Range("C20").Value = "Code changed 10010018"
$ws.Range("D20").Value = "gfndnvbx"
$ws.Range("E20").Value = "Code changed 10010018"

# Row 21: C21 changes from "dsgagass" to "Code changed 10010018"; D21 = "dsgagass"; E21 = "Code changed 10010018"
$ws.Range("C21").Value = "Code changed 10010018"
$ws.Range("D21").Value = "dsgagass"
$ws.Range("E21").Value = "Code changed 10010018"

# Row 22: C22 changes from "dsvsbsb" to "Code changed 10010018"; D22 = "dsvsbsb"; E22 = "Code changed 10010018"
$ws.Range("C22").Value = "Code changed 10010018"
$ws.Range("D22").Value = "dsvsbsb"
$ws.Range("E22").Value = "Code changed 10010018"

# Row 23: C23 changes from "nsngnsg" to "Code changed 10010018"; D23 = "nsngnsg"; E23 = "Code changed 10010018"
$ws.Range("C23").Value = "Code changed 10010018"
$ws.Range("D23").Value = "nsngnsg"
$ws.Range("E23").Value = "Code changed 10010018"

# New rows 24 and 25, all five columns = "Code changed 10010018"
$ws.Range("A24:E24").Value = "Code changed 10010018"
$ws.Range("A25:E25").Value = "Code changed 10010018"
